$d = $word.ActiveDocument

# The paragraph (styled "First Paragraph", in section "2 Задание") currently
# reads as one paragraph built from 5 runs:
#   "1.Реалтзация циклов в NASM." / " " / "2.Обработка аргументов командной
#   строки." / " " / "3.Задание для самостоятельной работы."
# Turn each of the two standalone-space runs into a paragraph break, so the
# single paragraph becomes three. Because "First Paragraph"'s style-for-the-
# following-paragraph is "Body Text", that's what Word applies to a
# paragraph newly created by pressing Enter here - so the two new
# paragraphs get the "Body Text" style explicitly.

$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text.StartsWith("1.Реалтзация")) {
        $target = $cand
        break
    }
}

# --- split #1: right before "2." ---
$pStart = $target.Range.Start
$sepPos = $target.Range.Text.IndexOf(" 2.")
$space1 = $d.Range($pStart + $sepPos, $pStart + $sepPos + 1)
$space1.InsertParagraphAfter()
$space1.Text = ""

$second = $target.Next()
$second.Style = "Body Text"

# --- split #2: right before "3." ---
$secondStart = $second.Range.Start
$sepPos2 = $second.Range.Text.IndexOf(" 3.")
$space2 = $d.Range($secondStart + $sepPos2, $secondStart + $sepPos2 + 1)
$space2.InsertParagraphAfter()
$space2.Text = ""

$third = $second.Next()
$third.Style = "Body Text"

Write-Output "paragraphs split"
